# Insert a new weekly price record at row 26, pushing the existing
# records (rows 26-79) down by one row (they become rows 27-80).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 26..79 down to 27..80 by inserting a new row at 26.
$ws.Rows(26).Insert()

# Fill in the new record's data in the now-empty row 26.
$ws.Cells.Item(26, 1).Value = 7
$ws.Cells.Item(26, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(26, 3).Value = "Ñuble"
$ws.Cells.Item(26, 4).Value = 45070
$ws.Cells.Item(26, 5).Value = 16
$ws.Cells.Item(26, 6).Value = 100112001
$ws.Cells.Item(26, 7).Value = "Berenjena"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 50
$ws.Cells.Item(26, 11).Value = 7000
$ws.Cells.Item(26, 12).Value = 7000
$ws.Cells.Item(26, 13).Value = 7000
$ws.Cells.Item(26, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(26, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(26, 16).Value = 117
$ws.Cells.Item(26, 17).Value = 60
$ws.Cells.Item(26, 18).Value = "Hortaliza"
